$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 662 (pushes the existing 662.. block down to 664..710),
# mirroring the weekly "Fruta / hortaliza" refresh described in the commit message.
$ws.Rows.Item(662).Insert()
$ws.Rows.Item(662).Insert()

# Row 662: new "Pintón" entry
$ws.Range("A662").Value = 7
$ws.Range("B662").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C662").Value = "Ñuble"
$ws.Range("D662").Value = 44826
$ws.Range("E662").Value = 16
$ws.Range("F662").Value = "Fruta"
$ws.Range("G662").Value = 100108
$ws.Range("H662").Value = "Tropicales y subtropicales"
$ws.Range("I662").Value = 100108006
$ws.Range("J662").Value = "Plátano"
$ws.Range("K662").Value = "Sin especificar"
$ws.Range("L662").Value = "Pintón"
$ws.Range("M662").Value = 80
$ws.Range("N662").Value = 21000
$ws.Range("O662").Value = 21000
$ws.Range("P662").Value = 21000
$ws.Range("Q662").Value = "$/caja 20 kilos"
$ws.Range("R662").Value = "Ecuador"
$ws.Range("S662").Value = 1050
$ws.Range("T662").Value = 20

# Row 663: new "Primera Pintón" entry
$ws.Range("A663").Value = 7
$ws.Range("B663").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C663").Value = "Ñuble"
$ws.Range("D663").Value = 44826
$ws.Range("E663").Value = 16
$ws.Range("F663").Value = "Fruta"
$ws.Range("G663").Value = 100108
$ws.Range("H663").Value = "Tropicales y subtropicales"
$ws.Range("I663").Value = 100108006
$ws.Range("J663").Value = "Plátano"
$ws.Range("K663").Value = "Sin especificar"
$ws.Range("L663").Value = "Primera Pintón"
$ws.Range("M663").Value = 160
$ws.Range("N663").Value = 22000
$ws.Range("O663").Value = 23000
$ws.Range("P663").Value = 22500
$ws.Range("Q663").Value = "$/caja 20 kilos"
$ws.Range("R663").Value = "Ecuador"
$ws.Range("S663").Value = 1125
$ws.Range("T663").Value = 20
